$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "formTestData"

# Fill in Gender and MobileNumber for row 2
$ws.Range("C2").Value = "Male"

# MobileNumber needs to be stored as text (leading zero must be kept)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0748146380"

# Update the active selection
$ws.Range("D5").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
